# Caelum results on BP4D.
# Updates the BP4D sheet (sheet1) with refreshed F1/Pr/Rec numbers for the
# "DISFA training with generic PCA and new alignment (masked and bigger),
# dynamic model" (row 8) and "SEMAINE trained (static), generic PCA" (row 9)
# result rows, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BP4D")

# --- Row 8 updates -------------------------------------------------------
$ws.Range("B8").Value = 0.3589
$ws.Range("C8").Value = 0.4729
$ws.Range("D8").Value = 0.4081
$ws.Range("E8").Value = 0.2986
$ws.Range("F8").Value = 0.4189
$ws.Range("G8").Value = 0.3487

$ws.Range("Q8").Value = 0.7231
$ws.Range("R8").Value = 0.7741
$ws.Range("S8").Value = 0.7477
$ws.Range("T8").Value = 0.8252
$ws.Range("U8").Value = 0.7614
$ws.Range("V8").Value = 0.792

# --- Row 9 updates -------------------------------------------------------
$ws.Range("Q9").Value = 0.7634
$ws.Range("R9").Value = 0.707
$ws.Range("S9").Value = 0.7341
$ws.Range("T9").Value = 0.859
$ws.Range("U9").Value = 0.7732
$ws.Range("V9").Value = 0.8139

# Previously-blank cells that now carry data (DISFA dynamic, p. video section)
$ws.Range("Z9").Value = 0.4131
$ws.Range("AA9").Value = 0.353
$ws.Range("AB9").Value = 0.3807
$ws.Range("AC9").Value = 0.5123
$ws.Range("AD9").Value = 0.8113
$ws.Range("AE9").Value = 0.628
$ws.Range("AF9").Value = 0.4255
$ws.Range("AG9").Value = 0.3859
$ws.Range("AH9").Value = 0.4047

# --- Selection -------------------------------------------------------------
$ws.Range("N9").Select()
